# "cambios en el modelado de la óptica"
# Re-models the provider/laboratory rows and replaces the "Montura" /
# "Montura con lentes" entries with "Mostrador" / "Bodega", widens column B
# and left-aligns its (now longer) descriptions, shrinks rows 18-19 back to
# a normal height, and drops the trailing blank rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6-8 group: "Proveedor lentes" -> "Proveedor" ---------------------
$ws.Range("A6").Value2 = "Proveedor"
$ws.Range("B6").Value2 = "Empresa a la que se le compran las monturas"
$ws.Range("D8").Value2 = "Monturas"

# --- Row 9-11 group: "Proveedor monturas" -> "Laboratorio" -----------------
$ws.Range("A9").Value2 = "Laboratorio"
$ws.Range("B9").Value2 = "Empresa a la que se le encargan las monturas junto con los lentes formulados o con determinada modificación. También se le encargan los arreglos de las mismas."
$ws.Range("D11").Value2 = "Monturas, Lentes, Arreglo"

# --- Rows 12-17 (Cliente / Paciente groups) are unchanged ------------------

# --- Row 18: "Montura" -> "Mostrador" --------------------------------------
$ws.Range("A18").Value2 = "Mostrador"
$ws.Range("B18").Value2 = "Estanterías donde se disponen las monturas disponibles a los clientes"
$ws.Range("D18").Value2 = "Monturas"
$ws.Rows.Item(18).RowHeight = 30
$ws.Range("A18").HorizontalAlignment = -4108

# --- Row 19: "Montura con lentes" -> "Bodega" ------------------------------
$ws.Range("A19").Value2 = "Bodega"
$ws.Range("B19").Value2 = "Lugar físico donde se guardan las reservas de monturas que no tienen en el mostrador."
$ws.Range("D19").Value2 = "Monturas, Lentes"
$ws.Rows.Item(19).RowHeight = 30
$ws.Range("A19").HorizontalAlignment = -4108

# --- Column B: wider, left aligned wrapped descriptions --------------------
$ws.Columns.Item(2).ColumnWidth = 49.285714285714285
$ws.Range("B6:B17").HorizontalAlignment = -4131

# --- Drop the trailing blank rows 20-24 ------------------------------------
$ws.Range("A20:E24").EntireRow.Delete() | Out-Null

# --- Match the saved selection ---------------------------------------------
$ws.Range("D20").Select() | Out-Null
